# Update res_bus vm_pu values for the 380 kV case (B1=1.02) across rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 2 = 1.02; 3 = 1.027286492920964; 4 = 1.028635010049834; 5 = 1.027392551042137; 6 = 1.037332505886597; 9 = 1.033220159888673; 10 = 1.032445210090636; 11 = 1.031451254316876; 12 = 1.030212412616719; 13 = 1.040123691758455; 14 = 1.014794193185123 }
    3 = @{ 2 = 1.02; 3 = 1.028445527965413; 4 = 1.029647051331626; 5 = 1.028381784931085; 6 = 1.038676180860914; 9 = 1.033476123635515; 10 = 1.033243289165005; 11 = 1.032270722215918; 12 = 1.031008871297636; 13 = 1.041275739484115; 14 = 1.015061666216107 }
    4 = @{ 2 = 1.02; 3 = 1.029194878085899; 4 = 1.030301635803313; 5 = 1.029021718869997; 6 = 1.039544959519403; 9 = 1.033639534863635; 10 = 1.033758577818463; 11 = 1.032800097011973; 12 = 1.031523459970915; 13 = 1.042019992310604; 14 = 1.015234267369492 }
    5 = @{ 2 = 1.02; 3 = 1.029509757682725; 4 = 1.030576758588191; 5 = 1.02929070815989; 6 = 1.039910036873449; 9 = 1.033707703116965; 10 = 1.033974937649984; 11 = 1.033022437526981; 12 = 1.031739608829578; 13 = 1.042332591591925; 14 = 1.015306716242162 }
    6 = @{ 2 = 1.02; 3 = 1.029562618778009; 4 = 1.03062294916071; 5 = 1.029335870388776; 6 = 1.039971325834876; 9 = 1.03371911780623; 10 = 1.034011249760912; 11 = 1.033059757271426; 12 = 1.031775890397925; 13 = 1.042385061768379; 14 = 1.015318874129156 }
    7 = @{ 2 = 1.02; 3 = 1.02919908609842; 4 = 1.030305312259693; 5 = 1.029025313269915; 6 = 1.039549838314904; 9 = 1.033640447812665; 10 = 1.03376146987899; 11 = 1.03280306875544; 12 = 1.031526348885274; 13 = 1.042024170392786; 14 = 1.015235235877162 }
    8 = @{ 2 = 1.02; 3 = 1.027678324436274; 4 = 1.028977091489018; 5 = 1.027726902312778; 6 = 1.037786747410742; 9 = 1.033307122739522; 10 = 1.032715157583938; 11 = 1.03172837926306; 12 = 1.030481740242936; 13 = 1.040513281443928; 14 = 1.014884684812043 }
    9 = @{ 2 = 1.02; 3 = 1.024993687293389; 4 = 1.026634451039057; 5 = 1.025437618188876; 6 = 1.034674700337787; 9 = 1.032702786655087; 10 = 1.030862780993875; 11 = 1.029827886206333; 12 = 1.028635041227412; 13 = 1.037841607483024; 14 = 1.014263342271197 }
    10 = @{ 2 = 1.02; 3 = 1.023200533607043; 4 = 1.025071177687319; 5 = 1.023910479579008; 6 = 1.03259627997779; 9 = 1.032288460111656; 10 = 1.029621983049885; 11 = 1.028556283547366; 12 = 1.027399836068135; 13 = 1.036054083885344; 14 = 1.013846655525646 }
    11 = @{ 2 = 1.02; 3 = 1.022423243550061; 4 = 1.024393887500639; 5 = 1.023248971540133; 6 = 1.03169537515836; 9 = 1.032106333898798; 10 = 1.029083292461378; 11 = 1.028004556930181; 12 = 1.026863997999271; 13 = 1.035278510375971; 14 = 1.013665637689646 }
    12 = @{ 2 = 1.02; 3 = 1.022134394169519; 4 = 1.024142252749447; 5 = 1.02300321949533; 6 = 1.031360594688066; 9 = 1.032038274855762; 10 = 1.028882984372318; 11 = 1.02779945221925; 12 = 1.026664814367688; 13 = 1.03499018962311; 14 = 1.013598310552879 }
    13 = @{ 2 = 1.02; 3 = 1.02219635923336; 4 = 1.024196231953619; 5 = 1.023055935908347; 6 = 1.031432412827225; 9 = 1.032052892281025; 10 = 1.028925960865321; 11 = 1.027843455541954; 12 = 1.026707546720754; 13 = 1.035052046235612; 14 = 1.013612756489412 }
    14 = @{ 2 = 1.02; 3 = 1.022399369822049; 4 = 1.024373088506622; 5 = 1.02322865839313; 6 = 1.031667705066272; 9 = 1.032100716469307; 10 = 1.029066739315345; 11 = 1.027987606364458; 12 = 1.026847536476158; 13 = 1.035254682571054; 14 = 1.013660074226804 }
    15 = @{ 2 = 1.02; 3 = 1.022524434140462; 4 = 1.024482047798577; 5 = 1.023335073271921; 6 = 1.031812657143491; 9 = 1.032130128281363; 10 = 1.029153449105488; 11 = 1.028076400021281; 12 = 1.026933768931447; 13 = 1.035379501841319; 14 = 1.013689216432515 }
    16 = @{ 2 = 1.02; 3 = 1.023252101808842; 4 = 1.025116119015289; 5 = 1.023954376366533; 6 = 1.032656050036777; 9 = 1.03230048986887; 10 = 1.029657704186122; 11 = 1.02859287622405; 12 = 1.027435376984043; 13 = 1.036105522916709; 14 = 1.013858656619858 }
    17 = @{ 2 = 1.02; 3 = 1.02370832104361; 4 = 1.025513751637975; 5 = 1.02434278173463; 6 = 1.033184835076482; 9 = 1.032406624530053; 10 = 1.029973629646041; 11 = 1.028916548742663; 12 = 1.027749757499651; 13 = 1.036560515744816; 14 = 1.013964783768126 }
    18 = @{ 2 = 1.02; 3 = 1.023974344879019; 4 = 1.025745647219669; 5 = 1.024569308285029; 6 = 1.03349317648202; 9 = 1.032468268674999; 10 = 1.030157767018764; 11 = 1.029105234036637; 12 = 1.027933035286808; 13 = 1.036825754682356; 14 = 1.014026629069475 }
    19 = @{ 2 = 1.02; 3 = 1.024065038451082; 4 = 1.025824711424667; 5 = 1.024646544001617; 6 = 1.033598297788693; 9 = 1.032489243244196; 10 = 1.030220529980258; 11 = 1.029169552671126; 12 = 1.027995512194732; 13 = 1.036916168844905; 14 = 1.014047707086575 }
    20 = @{ 2 = 1.02; 3 = 1.023659381448824; 4 = 1.025471093239257; 5 = 1.024301111979993; 6 = 1.033128110780326; 9 = 1.03239526443212; 10 = 1.029939747975616; 11 = 1.028881832866039; 12 = 1.027716037287619; 13 = 1.036511714930794; 14 = 1.013953403215257 }
    21 = @{ 2 = 1.02; 3 = 1.022339591852671; 4 = 1.024321010301566; 5 = 1.023177797000657; 6 = 1.031598421380783; 9 = 1.032086644740225; 10 = 1.02902528950215; 11 = 1.0279451622139; 12 = 1.026806317120724; 13 = 1.035195017823042; 14 = 1.013646142794458 }
    22 = @{ 2 = 1.02; 3 = 1.021509038661657; 4 = 1.023597565263276; 5 = 1.022471300801067; 6 = 1.030635808113041; 9 = 1.03189023514449; 10 = 1.02844909182381; 11 = 1.027355262372921; 12 = 1.026233474212772; 13 = 1.034365778205699; 14 = 1.013452440513609 }
    23 = @{ 2 = 1.02; 3 = 1.02194940251485; 4 = 1.02398111012996; 5 = 1.022845849461298; 6 = 1.031146188224175; 9 = 1.031994580262494; 10 = 1.028754663300966; 11 = 1.027668072622304; 12 = 1.026537231598505; 13 = 1.034805505704341; 14 = 1.0135551747699 }
    24 = @{ 2 = 1.02; 3 = 1.023681495395776; 4 = 1.025490368847204; 5 = 1.024319940821861; 6 = 1.033153742325933; 9 = 1.032400398382213; 10 = 1.029955058066172; 11 = 1.028897519808278; 12 = 1.027731274293698; 13 = 1.03653376638353; 14 = 1.013958545772965 }
    25 = @{ 2 = 1.02; 3 = 1.025688319438913; 4 = 1.02724034121827; 5 = 1.026029615640851; 6 = 1.035479881270094; 9 = 1.032861035932169; 10 = 1.031342695930946; 11 = 1.030320016313921; 12 = 1.029113170404962; 13 = 1.038533416935838; 14 = 1.014424406038485 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item($r, $c).Value = $data[$r][$c]
    }
}
